$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values are stored as literal text (e.g. "1.000", "29.352.44") in the
# source data. Force text number format before assignment so Excel does not
# auto-convert them to numeric values (which would drop formatting such as
# trailing zeros), then restore the default style so no extra formatting is
# left behind on the cell.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.352.44'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.55%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.874.06'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.73%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7136'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.77%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '241.79'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.38%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3112'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07769'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.55%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.12'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08440'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.879.14'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.45%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.244'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.7122'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.39%  '
$ws.Range('E15').Value = '  -0.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.357.90'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.54%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.079'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.62%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008238'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '240.46'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.21'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.97%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.123.35'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.62%  '
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.768'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.29%  '
$ws.Range('E24').Value = '  +0.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1592'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.13'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.048'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.44%  '
$ws.Range('E28').Value = '  +0.44%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.510'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.82%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.422'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.53%  '
$ws.Range('E31').Value = '  +2.61%  '
$ws.Range('E32').Value = '  -2.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05288'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.934'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.50%  '
$ws.Range('E35').Value = '  +1.31%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7426'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -7.89%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.701'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.78%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01874'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.85%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.226.71'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.07%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.734'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.32%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.529'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.61%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '110.60'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +8.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8884'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '72.91'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.000'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.020.69'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.802'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.78%  '
$ws.Range('E48').Value = '  +0.72%  '
$ws.Range('E49').Value = '  +2.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.424'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.82%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4315'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.10%  '
